# Updates cryptos list values (price & 1h volume change) per latest scrape.
# Rows 37-40 additionally swap which coin occupies the row (Coin/Link/Price/Volume).
# Price column (D) is stored as text in the sheet, so values are entered with a
# leading apostrophe to stop Excel re-interpreting dotted/zero-padded numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.783.11"
$ws.Range("E2").Value = "  -2.58%  "

$ws.Range("D3").Value = "'1.781.09"

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'310.56"
$ws.Range("E5").Value = "  -2.06%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").Value = "'0.5136"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").Value = "'0.3790"
$ws.Range("E8").Value = "  -2.30%  "

$ws.Range("D9").Value = "'0.07781"
$ws.Range("E9").Value = "  -7.61%  "

$ws.Range("D10").Value = "'41.29"
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").Value = "'1.084"
$ws.Range("E11").Value = "  -2.32%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").Value = "'6.204"
$ws.Range("E13").Value = "  -3.58%  "

$ws.Range("D14").Value = "'20.11"
$ws.Range("E14").Value = "  -4.20%  "

$ws.Range("D15").Value = "'1.782.87"
$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").Value = "'7.161"
$ws.Range("E16").Value = "  -4.71%  "

$ws.Range("D17").Value = "'91.45"
$ws.Range("E17").Value = "  -1.47%  "

$ws.Range("D18").Value = "'0.00001072"
$ws.Range("E18").Value = "  -5.60%  "

$ws.Range("D19").Value = "'0.06556"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("D21").Value = "'17.00"
$ws.Range("E21").Value = "  -4.23%  "

$ws.Range("D22").Value = "'5.911"
$ws.Range("E22").Value = "  -2.81%  "

$ws.Range("D23").Value = "'27.827.30"
$ws.Range("E23").Value = "  -2.56%  "

$ws.Range("D24").Value = "'10.99"
$ws.Range("E24").Value = "  -3.48%  "

$ws.Range("E25").Value = "  -1.66%  "

$ws.Range("D26").Value = "'159.29"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "'20.22"
$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("D28").Value = "'1.986.23"
$ws.Range("E28").Value = "  -2.22%  "

$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("D30").Value = "'124.95"
$ws.Range("E30").Value = "  -0.86%  "

$ws.Range("D31").Value = "'0.1072"
$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").Value = "'1.030"
$ws.Range("E32").Value = "  -5.85%  "

$ws.Range("D33").Value = "'3.632"
$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").Value = "'5.473"
$ws.Range("E34").Value = "  -4.67%  "

$ws.Range("D35").Value = "'0.07064"
$ws.Range("E35").Value = "  -5.93%  "

$ws.Range("D36").Value = "'0.02316"
$ws.Range("E36").Value = "  -2.09%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'8.701"
$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2122"
$ws.Range("E38").Value = "  -4.89%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'11.52"
$ws.Range("E39").Value = "  +2.45%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'5.000"
$ws.Range("E40").Value = "  -3.69%  "

$ws.Range("D41").Value = "'0.6082"
$ws.Range("E41").Value = "  -3.81%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "'1.152"
$ws.Range("E43").Value = "  -3.58%  "

$ws.Range("D44").Value = "'1.320"
$ws.Range("E44").Value = "  -5.90%  "

$ws.Range("D45").Value = "'13.12"
$ws.Range("E45").Value = "  -2.81%  "

$ws.Range("D46").Value = "'0.5919"
$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("D47").Value = "'3.710"
$ws.Range("E47").Value = "  -1.81%  "

$ws.Range("D48").Value = "'127.87"
$ws.Range("E48").Value = "  +1.48%  "

$ws.Range("D49").Value = "'1.204"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("D50").Value = "'1.893"
$ws.Range("E50").Value = "  -4.85%  "

$ws.Range("D51").Value = "'0.06790"
$ws.Range("E51").Value = "  -2.69%  "
